# Weekly update: insert a new price record as row 85 for
# "Hortaliza, Agrícola del Norte S.A. de Arica - Cebollín baby", pushing the
# existing rows 85-110 down to 86-111 (dimension grows from A1:R110 to A1:R111).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 85; everything below shifts down
# by one (old row 85 -> 86, ..., old row 110 -> 111).
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with the new weekly record. All
# non-numeric/template fields (market, region, product, quality, unit,
# origin, classification, etc.) are identical across every row in this
# sheet.
$ws.Range("A85").Value = 1
$ws.Range("B85").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C85").Value = "Arica y Parinacota"
$ws.Range("D85").Value = 44876
$ws.Range("E85").Value = 15
$ws.Range("F85").Value = 100112038
$ws.Range("G85").Value = "Cebollín baby"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 200
$ws.Range("K85").Value = 1000
$ws.Range("L85").Value = 1200
$ws.Range("M85").Value = 1100
$ws.Range("N85").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O85").Value = "Región de Arica y Parinacota"
$ws.Range("P85").Value = 550
$ws.Range("Q85").Value = 2
$ws.Range("R85").Value = "Hortaliza"
